# Generate Report for Handback
# Adds a new handback record (fe721b5d-0fa7-434f-8789-e65e1e9fafb1) as row 4
# to the Overview, zh-cn and de-de sheets/tables.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ovWs = $wb.Worksheets.Item("Overview")
$ovTable = $ovWs.ListObjects.Item("Overview")
$ovTable.ListRows.Add() | Out-Null

$ovWs.Range("A4").Value = "fe721b5d-0fa7-434f-8789-e65e1e9fafb1.md"
$ovWs.Range("B4").Value = "e2e\fe721b5d-0fa7-434f-8789-e65e1e9fafb1.md"
$ovWs.Range("C4").Value = ".md"
$ovWs.Range("E4").Value = "Handed back: in sync with en-US"
$ovWs.Range("F4").Value = "Handed back: in sync with en-US"
$ovWs.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ovWs.Range("G4").Value = "2016-08-21 00:51:30"

$ovWs.Hyperlinks.Add(
    $ovWs.Range("B4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f130e10cab627b8e564c031aedd92e310d428cb2/e2e/fe721b5d-0fa7-434f-8789-e65e1e9fafb1.md",
    "",
    "",
    "e2e\fe721b5d-0fa7-434f-8789-e65e1e9fafb1.md"
) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$zhWs = $wb.Worksheets.Item("zh-cn")
$zhTable = $zhWs.ListObjects.Item("zh-cn")
$zhTable.ListRows.Add() | Out-Null

$zhWs.Range("A4").Value = "fe721b5d-0fa7-434f-8789-e65e1e9fafb1.md"
$zhWs.Range("B4").Value = ".md"
$zhWs.Range("C4").Value = "Handed back: in sync with en-US"
$zhWs.Range("D4").Value = "e2e"
$zhWs.Range("E4").Value = "ht"
$zhWs.Range("F4").Value = "'True"
$zhWs.Range("G4").Value = "fe721b5d-0fa7-434f-8789-e65e1e9fafb1.ee0f9686e6d1b7fa1b6ce94b379570f14e04ced2.zh-cn.xlf"
$zhWs.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zhWs.Range("H4").Value = "2016-08-21 00:51:26"
$zhWs.Range("I4").Value = "fe721b5d-0fa7-434f-8789-e65e1e9fafb1.md"
$zhWs.Range("J4").Value = "fe721b5d-0fa7-434f-8789-e65e1e9fafb1.ee0f9686e6d1b7fa1b6ce94b379570f14e04ced2.zh-cn.xlf"
$zhWs.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zhWs.Range("K4").Value = "2016-08-21 00:51:43"
$zhWs.Range("L4").Value = "'"
$zhWs.Range("M4").Value = "'True"
$zhWs.Range("N4").Value = "'"
$zhWs.Range("O4").Value = "'False"
$zhWs.Range("P4").Value = "'"

$zhWs.Hyperlinks.Add(
    $zhWs.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f130e10cab627b8e564c031aedd92e310d428cb2/e2e/fe721b5d-0fa7-434f-8789-e65e1e9fafb1.md",
    "",
    "",
    "fe721b5d-0fa7-434f-8789-e65e1e9fafb1.md"
) | Out-Null
$zhWs.Hyperlinks.Add(
    $zhWs.Range("I4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/98bf03374ea8087a340fefc44c71bb9544d46002/e2e/fe721b5d-0fa7-434f-8789-e65e1e9fafb1.md",
    "",
    "",
    "fe721b5d-0fa7-434f-8789-e65e1e9fafb1.md"
) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$deWs = $wb.Worksheets.Item("de-de")
$deTable = $deWs.ListObjects.Item("de-de")
$deTable.ListRows.Add() | Out-Null

$deWs.Range("A4").Value = "fe721b5d-0fa7-434f-8789-e65e1e9fafb1.md"
$deWs.Range("B4").Value = ".md"
$deWs.Range("C4").Value = "Handed back: in sync with en-US"
$deWs.Range("D4").Value = "e2e"
$deWs.Range("E4").Value = "ht"
$deWs.Range("F4").Value = "'True"
$deWs.Range("G4").Value = "fe721b5d-0fa7-434f-8789-e65e1e9fafb1.ee0f9686e6d1b7fa1b6ce94b379570f14e04ced2.de-de.xlf"
$deWs.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$deWs.Range("H4").Value = "2016-08-21 00:51:30"
$deWs.Range("I4").Value = "fe721b5d-0fa7-434f-8789-e65e1e9fafb1.md"
$deWs.Range("J4").Value = "fe721b5d-0fa7-434f-8789-e65e1e9fafb1.ee0f9686e6d1b7fa1b6ce94b379570f14e04ced2.de-de.xlf"
$deWs.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$deWs.Range("K4").Value = "2016-08-21 00:51:49"
$deWs.Range("L4").Value = "'"
$deWs.Range("M4").Value = "'True"
$deWs.Range("N4").Value = "'"
$deWs.Range("O4").Value = "'False"
$deWs.Range("P4").Value = "'"

$deWs.Hyperlinks.Add(
    $deWs.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f130e10cab627b8e564c031aedd92e310d428cb2/e2e/fe721b5d-0fa7-434f-8789-e65e1e9fafb1.md",
    "",
    "",
    "fe721b5d-0fa7-434f-8789-e65e1e9fafb1.md"
) | Out-Null
$deWs.Hyperlinks.Add(
    $deWs.Range("I4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/95e6a97b90e1c0ebf33b6ec65b15d5d1ce2e139b/e2e/fe721b5d-0fa7-434f-8789-e65e1e9fafb1.md",
    "",
    "",
    "fe721b5d-0fa7-434f-8789-e65e1e9fafb1.md"
) | Out-Null

Write-Host "Handback row appended to Overview, zh-cn and de-de sheets."
